# Generate Report for Handoff
# Source file's GUID-based name (and its dependent .xlf report names) changed
# because a new handoff was generated. Update the cell text, the matching
# hyperlink display text (keeping the same target address/relationship) and
# the handoff timestamps for both locales.

$wb = $excel.ActiveWorkbook

$oldId = "3f9a74ae-5341-442d-91e0-004b45e2938d"
$newId = "8177b148-a9e4-41f9-b46f-e2e11509f714"

$oldHash = "f60c3d65f0b51dd1ee9a900201a81c947de30a2e"
$newHash = "c088dfbaf72c5aa9698b7d3d0fc4b025db8f02d7"

$oldMd = "$oldId.md"
$newMd = "$newId.md"

$oldZhXlf = "$oldId.$oldHash.zh-cn.xlf"
$newZhXlf = "$newId.$newHash.zh-cn.xlf"

$oldDeXlf = "$oldId.$oldHash.de-de.xlf"
$newDeXlf = "$newId.$newHash.de-de.xlf"

$oldZhTime = "2016-03-08 21:13:04"
$newZhTime = "2016-03-08 21:13:47"

$oldDeTime = "2016-03-08 21:13:11"
$newDeTime = "2016-03-08 21:13:55"

function Update-HyperlinkDisplay($ws, $oldText, $newText) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.TextToDisplay -eq $oldText) {
            $hl.TextToDisplay = $newText
        }
    }
}

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMd
Update-HyperlinkDisplay $wsOverview $oldMd $newMd

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = $newMd
$wsZh.Range("C2").Value = $newZhXlf
$wsZh.Range("D2").Value = $newZhTime
Update-HyperlinkDisplay $wsZh $oldMd $newMd
Update-HyperlinkDisplay $wsZh $oldZhXlf $newZhXlf

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = $newMd
$wsDe.Range("C2").Value = $newDeXlf
$wsDe.Range("D2").Value = $newDeTime
Update-HyperlinkDisplay $wsDe $oldMd $newMd
Update-HyperlinkDisplay $wsDe $oldDeXlf $newDeXlf
